$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook previously contained results for GKD, MDG and SOM instances.
# This update keeps GKD instances only: the MDG-a/MDG-b/SOM-a/SOM-a/SOM-b
# rows (rows 26-30, instance names) are cleared out, which also drops their
# shared-string entries and lets the AVERAGE() summary formulas in rows 3-4
# recompute over the smaller (GKD-only) data range.
$ws.Range("A26:O30").ClearContents() | Out-Null

# Restore the current selection to match the saved workbook view.
$ws.Range("T24").Select() | Out-Null
